$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style changes ---
# The title font (row 1) loses its explicit 14pt size and the header font
# (row 2) is merged with it: both become a single bold, white font (used
# on the dark-blue header fill, and now also on the title cell).
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2:K2").Font.Color = 16777215
$ws.Range("A2:K2").Font.Bold = $true

# --- Data changes: "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) refreshed ---
$ws.Range("H3").Value = 113
$ws.Range("I3").Value = "16-Sep-2025"

$ws.Range("H4").Value = 680
$ws.Range("I4").Value = "16-Sep-2025"

$ws.Range("H5").Value = 679
$ws.Range("I5").Value = "16-Sep-2025"

$ws.Range("H6").Value = 680
$ws.Range("I6").Value = "16-Sep-2025"

$ws.Range("H7").Value = 679
$ws.Range("I7").Value = "16-Sep-2025"
